$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 8, 9, 10 (old Inflammatory-Mac target rows, and extra MuSCs-sender rows no longer present)
$ws.Range("A8:T10").EntireRow.Delete() | Out-Null

# Updated data for rows 2-7 with new TPM-derived values
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..J ligand stats, K..P receptor stats, Q..T edge stats

$data = @(
    @("ECs",  "Vtn", "Tnfrsf11b", "FAPs",  2, 0.6666666666666666, 1.440985666666667, 4.322957000000001, 0.1098365531732288, 0.1230162332390494, 3, 1, 2.311298, 6.933894, 0.9776743782722859, 0.9850046587426607, 3.330547289395335, 29.97492560455801, 0.1073843838352073, 0.1211715628414374),
    @("ECs",  "Vtn", "Tnfrsf11b", "MuSCs", 2, 0.6666666666666666, 1.440985666666667, 4.322957000000001, 0.1098365531732288, 0.1230162332390494, 1, 0.5, 0.0527795, 0.105559, 0.02232562172771408, 0.01499534125733917, 0.07605450299383335, 0.4563270179630001, 0.00245216933802146, 0.001844670397611975),
    @("FAPs", "Vtn", "Tnfrsf11b", "FAPs",  3, 1, 7.461641333333333, 22.384924, 0.5687502547919595, 0.6369966279614609, 3, 1, 2.311298, 6.933894, 0.9776743782722859, 0.9850046587426607, 17.24607669045067, 155.214690214056, 0.5560525517459332, 0.6274446461454045),
    @("FAPs", "Vtn", "Tnfrsf11b", "MuSCs", 3, 1, 7.461641333333333, 22.384924, 0.5687502547919595, 0.6369966279614609, 1, 0.5, 0.0527795, 0.105559, 0.02232562172771408, 0.01499534125733917, 0.3938216987526666, 2.362930192516, 0.01269770304602629, 0.009551981816056425),
    @("MuSCs","Vtn", "Tnfrsf11b", "FAPs",  2, 1, 4.2167365, 8.433472999999999, 0.3214131920348118, 0.2399871387994896, 3, 1, 2.311298, 6.933894, 0.9776743782722859, 0.9850046587426607, 9.746134638977001, 58.476807833862, 0.3142374426911455, 0.2363884497558188),
    @("MuSCs","Vtn", "Tnfrsf11b", "MuSCs", 2, 1, 4.2167365, 8.433472999999999, 0.3214131920348118, 0.2399871387994896, 1, 0.5, 0.0527795, 0.105559, 0.02232562172771408, 0.01499534125733917, 0.22255724410175, 0.8902289764069999, 0.007175749343666333, 0.003598689043670768)
)

$r = 2
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}
